$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column E header "week", matching the style of the other headers (D1) ---
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value = "week"

# --- Column A: updated date-serial values (A2:A22) ---
$aVals = @(45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45012,45011,45012,45013)
$aArr = New-Object 'object[,]' 21,1
for ($i = 0; $i -lt 21; $i++) { $aArr[$i,0] = $aVals[$i] }
$ws.Range("A2:A22").Value = $aArr

# --- Column E: new "week" number values (E2:E22) ---
$eVals = @(13,13,13,13,13,13,13,13,13,13,13,13,13,13,13,13,13,13,12,13,13)
$eArr = New-Object 'object[,]' 21,1
for ($i = 0; $i -lt 21; $i++) { $eArr[$i,0] = $eVals[$i] }
$ws.Range("E2:E22").Value = $eArr

# --- Column C / D: updated check-out / work_hours text on specific rows ---
$ws.Range("C4").Value = "22:40:55"
$ws.Range("D4").Value = "11:25:17"

$ws.Range("C7").Value = "22:40:55"
$ws.Range("D7").Value = "11:18:29"

$ws.Range("C10").Value = "22:40:55"
$ws.Range("D10").Value = "11:12:27"

$ws.Range("C15").Value = "22:40:55"
$ws.Range("D15").Value = "11:01:29"

$ws.Range("C19").Value = "22:40:55"
$ws.Range("D19").Value = "10:50:51"

$ws.Range("C22").Value = "17:40:35"
$ws.Range("D22").Value = "05:27:06"
